# Hindalco prices: add the newest ready-reckoner entry (13.09.2025) at the
# top of the data table (row 2), pushing all existing rows down by one.
#
# Before: rows 2..31 hold Sl.no. 30..1 (newest first).
# After : a new row 2 holds Sl.no. 31 (13.09.2025 / 272), and the old
#         rows 2..31 become rows 3..32 (Sl.no. 30..1), unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$descText = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$gradeText = "P1020"

# 1. Insert a new blank row above row 2; existing rows 2-31 shift to 3-32.
$ws.Rows("2:2").Insert()

# 2. Copy the (now shifted) row 3 formatting onto the new row 2 so the new
#    row uses the same cell styles as every other data row.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new row with the newest circular's data.
$ws.Range("A2").Value = 31
$ws.Range("B2").Value = $descText
$ws.Range("C2").Value = $gradeText
$ws.Range("D2").Value = 272
$ws.Range("E2").Value = "13.09.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-september-2025.pdf"

# 4. Rebuild the hyperlinks for column F top to bottom (F2..F12) so the
#    relationship ids + targets line up with the new row order.
$ws.Range("F2").Hyperlinks.Delete()

$links = @(
    @("F2", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-september-2025.pdf"),
    @("F3", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-september-2025.pdf"),
    @("F4", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf"),
    @("F5", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf"),
    @("F6", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf"),
    @("F7", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf"),
    @("F8", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf"),
    @("F9", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf"),
    @("F10", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf"),
    @("F11", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf"),
    @("F12", "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link[0]), $link[1])
}

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" look
# (underline + theme colour); the source data keeps the plain data-row
# style, so restore it from a neighbouring (non-hyperlinked) cell.
$ws.Range("A2").Copy()
foreach ($link in $links) {
    $ws.Range($link[0]).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

Write-Output "done"
